# Append a new transaction row (row 3) to the Transactions sheet, matching
# the existing rows' convention of storing every value as literal text
# (dates and numbers included) rather than as native number/date types.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = "A3:E3"

# Force text storage for the new cells (so "2025-09-04", "10", "32" are not
# auto-coerced into a date serial / numeric types), then restore the default
# "Normal" style so no stray number-format style is left behind on the cells.
$ws.Range($newRow).NumberFormat = "@"

$ws.Range("A3").Value = "2025-09-04"
$ws.Range("B3").Value = "MEBL"
$ws.Range("C3").Value = "Buy"
$ws.Range("D3").Value = "10"
$ws.Range("E3").Value = "32"

$ws.Range($newRow).Style = "Normal"
